$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Modeller – " / bookmark / "SD, Package"  ->  single run "Modeller – SD, Package"
#    (drops the stray _GoBack bookmark that sat between the two runs)
# ---------------------------------------------------------------------------
$rngModeller = $d.Content
$rngModeller.Find.Execute("Modeller – SD, Package", $true, $false, $false, $false, $false, $true, 1, $false, "Modeller – SD, Package", 2)

# ---------------------------------------------------------------------------
# 2) Strike through the "User Story 8: Navbar" bullet
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "User Story 8*") {
        $p.Range.Font.StrikeThrough = 1
    }
}

# ---------------------------------------------------------------------------
# 3) "...ihht. Ny User Story" -> "...ihht. Nye User Storys og tasks"
#    Split into the same run layout as the authored edit, with a _GoBack
#    bookmark sitting between "Nye" and " User Story".
# ---------------------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Ny User Story*") {
        $target = $p
    }
}
$pStart = $target.Range.Start
$pEnd = $target.Range.End

# 3a) "Ny" -> "Ny" + new run "e"  (". Ny" | "e")
$r1 = $d.Range($pStart, $pEnd)
$r1.Find.Execute("Ny", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r1.Collapse(0)
$r1.InsertAfter("e")
# force the just-inserted text into its own run without leaving formatting residue
$d.Bookmarks.Add("zzsplit1", $r1)
$d.Bookmarks("zzsplit1").Delete()

# 3b) place the _GoBack bookmark right after "Nye", before " User Story"
$pEnd = $target.Range.End
$r2 = $d.Range($pStart, $pEnd)
$r2.Find.Execute("Nye", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r2.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r2)

# 3c) " User Story" -> " User Story" + new run "s og tasks"
$pEnd = $target.Range.End
$r3 = $d.Range($pStart, $pEnd)
$r3.Find.Execute("User Story", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r3.Collapse(0)
$r3.InsertAfter("s og tasks")
$d.Bookmarks.Add("zzsplit2", $r3)
$d.Bookmarks("zzsplit2").Delete()

# ---------------------------------------------------------------------------
# 4) Add a new bullet after "Sæt motiv på user storys":
#    "User story 10: Opret bruger på login-side" (struck through)
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$endRng = $lastPara.Range
$endRng.Collapse(0)
$endRng.InsertParagraphAfter()
$newPara = $d.Paragraphs($d.Paragraphs.Count)
$newPara.Range.Text = "User story 10: Opret bruger på login-side"
$newPara.Range.Font.StrikeThrough = 1

Write-Output "done"
